$wb = $excel.ActiveWorkbook

# "Tabla" sheet holds the underlying schedule data (section/day/start/end/room).
# "Computación Cuántica (Sección A)" was recorded as starting at 9:00, but it
# actually starts at 10:00 — fix the start time.
$wsTabla = $wb.Worksheets.Item("Tabla")
$wsTabla.Range("C5").Value = "10:00"

# "Horario" sheet is the rendered weekly grid view built from that data.
# Because the course really starts at 10:00 (not 9:00), the entry that was
# incorrectly shown in the 9:00 row under "Jueves" must be removed, leaving
# that cell blank (its 10:00-row entry, D3/E3, was already correct).
$wsHorario = $wb.Worksheets.Item("Horario")
$wsHorario.Range("E2").ClearContents()
